# Regenerate the handback-status report for the new CI run: the two
# source files tracked by this workbook were re-handed-off/-back under
# new GUIDs, and the zh-cn / de-de Xliff round-trip timestamps (plus the
# generated Xliff file names) were refreshed. Update every sheet that
# caches these values (Overview, zh-cn, de-de) including the cached
# hyperlink display text.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "1ca1ff6a-5ca0-47c2-8b7d-b99046b879ed"
$newGuid1 = "b87ae8ea-0e99-49ce-9582-2dcb932b0314"

$oldGuid2 = "9cb632ca-9933-4358-945b-2574f1e9564e"
$newGuid2 = "ffff072ee604-92e1-47b0-b8f1-2f921b68a245"

$newFile1 = "$newGuid1.md"
$newFile2 = "$newGuid2.md"

$newPath1 = "e2e\$newGuid1.md"
$newPath2 = "e2e\$newGuid2.md"

$newXlfStem = "$newGuid1.7a8ac1c9e2a5fc0e24b2b886b13582292981f067"
$newZhCnXlf = "$newXlfStem.zh-cn.xlf"
$newDeDeXlf = "$newXlfStem.de-de.xlf"

$newOverviewDate = "2016-08-23 07:01:11"
$newZhCnHandoffDate = "2016-08-23 07:00:58"
$newZhCnHandbackDate = "2016-08-23 07:01:32"
$newDeDeHandoffDate = "2016-08-23 07:01:11"
$newDeDeHandbackDate = "2016-08-23 07:01:40"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("B2").Value = $newPath1
$wsOverview.Range("G2").Value = $newOverviewDate

$wsOverview.Range("A3").Value = $newFile2
$wsOverview.Range("B3").Value = $newPath2
$wsOverview.Range("G3").Value = $newOverviewDate

# Re-create the hyperlinks so their cached "display" text (and the
# underlying cell's visible text) reflects the new file names, while the
# link targets (pointing at the historical git blob URLs) are preserved.
# (NOTE: collect addresses via foreach - indexed .Item(n) access on the
# Hyperlinks collection does not reliably resolve properties here.)
$ovAddrs = @()
foreach ($hl in $wsOverview.Hyperlinks) { $ovAddrs += $hl.Address }
$ovAddr1 = $ovAddrs[0]
$ovAddr2 = $ovAddrs[1]
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $ovAddr1, "", "", $newPath1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $ovAddr2, "", "", $newPath2)

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhAddrs = @()
foreach ($hl in $wsZhCn.Hyperlinks) { $zhAddrs += $hl.Address }
$zhAddrA2 = $zhAddrs[0]
$zhAddrI2 = $zhAddrs[1]
$zhAddrA3 = $zhAddrs[2]
$zhAddrI3 = $zhAddrs[3]

$wsZhCn.Range("A2").Value = $newFile1
$wsZhCn.Range("G2").Value = $newZhCnXlf
$wsZhCn.Range("H2").Value = $newZhCnHandoffDate
$wsZhCn.Range("I2").Value = $newFile1
$wsZhCn.Range("J2").Value = $newZhCnXlf
$wsZhCn.Range("K2").Value = $newZhCnHandbackDate

$wsZhCn.Range("A3").Value = $newFile2
$wsZhCn.Range("G3").Value = $newZhCnXlf
$wsZhCn.Range("H3").Value = $newZhCnHandoffDate
$wsZhCn.Range("I3").Value = $newFile2
$wsZhCn.Range("J3").Value = $newZhCnXlf
$wsZhCn.Range("K3").Value = $newZhCnHandbackDate

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhAddrA2, "", "", $newFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $zhAddrI2, "", "", $newFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $zhAddrA3, "", "", $newFile2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $zhAddrI3, "", "", $newFile2)

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deAddrs = @()
foreach ($hl in $wsDeDe.Hyperlinks) { $deAddrs += $hl.Address }
$deAddrA2 = $deAddrs[0]
$deAddrI2 = $deAddrs[1]
$deAddrA3 = $deAddrs[2]
$deAddrI3 = $deAddrs[3]

$wsDeDe.Range("A2").Value = $newFile1
$wsDeDe.Range("G2").Value = $newDeDeXlf
$wsDeDe.Range("H2").Value = $newDeDeHandoffDate
$wsDeDe.Range("I2").Value = $newFile1
$wsDeDe.Range("J2").Value = $newDeDeXlf
$wsDeDe.Range("K2").Value = $newDeDeHandbackDate

$wsDeDe.Range("A3").Value = $newFile2
$wsDeDe.Range("G3").Value = $newDeDeXlf
$wsDeDe.Range("H3").Value = $newDeDeHandoffDate
$wsDeDe.Range("I3").Value = $newFile2
$wsDeDe.Range("J3").Value = $newDeDeXlf
$wsDeDe.Range("K3").Value = $newDeDeHandbackDate

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $deAddrA2, "", "", $newFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $deAddrI2, "", "", $newFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $deAddrA3, "", "", $newFile2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $deAddrI3, "", "", $newFile2)
